# Add a new "Poster" URL column (D) with a link for each movie row,
# matching the commit "Added pictures and css".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "https://images-na.ssl-images-amazon.com/images/M/MV5BMTkxMTA5OTAzMl5BMl5BanBnXkFtZTgwNjA5MDc3NjE@._V1_SY1000_CR0,0,673,1000_AL_.jpg"
$ws.Range("D2").Value = "https://images-na.ssl-images-amazon.com/images/M/MV5BNzQzOTk3OTAtNDQ0Zi00ZTVkLWI0MTEtMDllZjNkYzNjNTc4L2ltYWdlXkEyXkFqcGdeQXVyNjU0OTQ0OTY@._V1_SY1000_CR0,0,665,1000_AL_.jpg"
$ws.Range("D3").Value = "https://images-na.ssl-images-amazon.com/images/M/MV5BMTc0NDQzNTA2Ml5BMl5BanBnXkFtZTcwNzI2OTQzMw@@._V1_.jpg"
$ws.Range("D4").Value = "http://www.gstatic.com/tv/thumb/movieposters/1587/p1587_p_v8_ag.jpg"
$ws.Range("D5").Value = "https://images-na.ssl-images-amazon.com/images/M/MV5BMjAxNjUyNjUwN15BMl5BanBnXkFtZTcwMDgwOTIyOA@@._V1_.jpg"

$ws.Range("D5").Select()
